$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark ambiguous numeric-looking D-column price cells as Text
# before assigning, so Excel keeps them as strings (matches original inlineStr type)
# instead of auto-converting to numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: write the updated cell values (order follows the diff)
$ws.Range("D2").Value = '29.591.66'
$ws.Range("D3").Value = '1.852.10'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '241.02'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").Value = '0.6305'
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.07489'
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").Value = '0.2912'
$ws.Range("D10").Value = '24.85'
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("D11").Value = '0.07749'
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '1.849.90'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '5.045'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").Value = '0.6828'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '0.00001038'
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").Value = '82.70'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '6.287'
$ws.Range("E17").Value = '  +2.96%  '
$ws.Range("D18").Value = '29.588.90'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '230.37'
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = '12.40'
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '7.577'
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '159.37'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = '8.527'
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("D26").Value = '0.1369'
$ws.Range("E26").Value = '  -1.66%  '
$ws.Range("D27").Value = '17.61'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").Value = '0.06675'
$ws.Range("E28").Value = '  +17.99%  '
$ws.Range("D29").Value = '1.441'
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").Value = '1.484'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").Value = '4.110'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '4.112'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '1.845'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").Value = '1.147'
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D35").Value = '0.7005'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").Value = '2.564'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = '0.01873'
$ws.Range("E37").Value = '  +2.57%  '
$ws.Range("D38").Value = '1.265.84'
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").Value = '2.844'
$ws.Range("E39").Value = '  +4.51%  '
$ws.Range("D40").Value = '6.817'
$ws.Range("E40").Value = '  +6.74%  '
$ws.Range("D41").Value = '0.9385'
$ws.Range("E41").Value = '  +4.43%  '
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '2.014.74'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").Value = '101.45'
$ws.Range("D45").Value = '66.41'
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000120'
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.735'
$ws.Range("E47").Value = '  +3.83%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.115'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1168'
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.004'
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.3958'
$ws.Range("E51").Value = '  -0.93%  '
